# Commit: "Fruta / hortaliza, semanal"
# This edit inserts two new weekly price-observation rows (rows 199 and 200)
# into the "Poroto granado" sheet, pushing all existing rows (old 199..280)
# down by two (new 201..282). The workbook's dimension grows from A1:R280
# to A1:R282 automatically as a result of the row insertion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 199, shifting everything
# below (old rows 199-280) down to rows 201-282.
$ws.Rows("199:200").Insert()

# ---- New row 199 ----
$ws.Cells.Item(199,1).Value  = 6
$ws.Cells.Item(199,2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(199,3).Value  = "Metropolitana"
$ws.Cells.Item(199,4).Value  = 44510
$ws.Cells.Item(199,5).Value  = 13
$ws.Cells.Item(199,6).Value  = 100112030
$ws.Cells.Item(199,7).Value  = "Poroto granado"
$ws.Cells.Item(199,8).Value  = "Sin especificar"
$ws.Cells.Item(199,9).Value  = "Primera"
$ws.Cells.Item(199,10).Value = 400
$ws.Cells.Item(199,11).Value = 30000
$ws.Cells.Item(199,12).Value = 32000
$ws.Cells.Item(199,13).Value = 30850
$ws.Cells.Item(199,14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(199,15).Value = "Perú"
$ws.Cells.Item(199,16).Value = 1234
$ws.Cells.Item(199,17).Value = 25
$ws.Cells.Item(199,18).Value = "Hortaliza"

# ---- New row 200 ----
$ws.Cells.Item(200,1).Value  = 6
$ws.Cells.Item(200,2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(200,3).Value  = "Metropolitana"
$ws.Cells.Item(200,4).Value  = 44510
$ws.Cells.Item(200,5).Value  = 13
$ws.Cells.Item(200,6).Value  = 100112030
$ws.Cells.Item(200,7).Value  = "Poroto granado"
$ws.Cells.Item(200,8).Value  = "Sin especificar"
$ws.Cells.Item(200,9).Value  = "Segunda"
$ws.Cells.Item(200,10).Value = 100
$ws.Cells.Item(200,11).Value = 25000
$ws.Cells.Item(200,12).Value = 25000
$ws.Cells.Item(200,13).Value = 25000
$ws.Cells.Item(200,14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(200,15).Value = "Perú"
$ws.Cells.Item(200,16).Value = 1000
$ws.Cells.Item(200,17).Value = 25
$ws.Cells.Item(200,18).Value = "Hortaliza"
